$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("4090845", "Prix Garantie Batterien AA/LR6 10 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/prix-garantie-batterien-aalr6-10-stueck/p/4090845", "10ST", 10, 3.5, "Coop", "4.75", "0.48/1ST", "Preis pro 1 Stück", "0.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Prix Garantie Batterien AA/LR6 10 Stück 4.75 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4090846", "Prix Garantie Batterien AAA/LR03 10 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/prix-garantie-batterien-aaalr03-10-stueck/p/4090846", "10ST", 13, 4, "Coop", "4.70", "0.47/1ST", "Preis pro 1 Stück", "0.47", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Prix Garantie Batterien AAA/LR03 10 Stück 4.70 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5920020", "BASIC Alkaline Batterien LR6/AA, 12 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/basic-alkaline-batterien-lr6aa-12-stueck/p/5920020", "12ST", $null, 0, "Coop", "9.95", "0.83/1ST", "Preis pro 1 Stück", "0.83", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "BASIC Alkaline Batterien LR6/AA, 12 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5920351", "BASIC Alkaline Batterien LR03/AAA 12 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/basic-alkaline-batterien-lr03aaa-12-stueck/p/5920351", "12ST", 3, 4, "Coop", "9.95", "0.83/1ST", "Preis pro 1 Stück", "0.83", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "BASIC Alkaline Batterien LR03/AAA 12 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4488074", "Varta Knopfzellen CR2032 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2032-2-stueck/p/4488074", "2ST", 4, 4.5, "Varta", "8.95", "4.48/1ST", "Preis pro 1 Stück", "4.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Knopfzellen CR2032 2 Stück 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6762335", "Duracell Batterien Optimum AA/LR6 8 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-8-stueck/p/6762335", "8ST", 1, 5, "Duracell", "19.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Duracell Batterien Optimum AA/LR6 8 Stück 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4014528", "Varta Longlife Batterien AAA/LR03 10 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-batterien-aaalr03-10-stueck/p/4014528", "10ST", 1, 5, "Varta", "16.95", "1.70/1ST", "Preis pro 1 Stück", "1.70", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Varta Longlife Batterien AAA/LR03 10 Stück 16.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6762336", "Duracell Batterien Optimum AAA/LR03 8 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-optimum-aaalr03-8-stueck/p/6762336", "8ST", $null, 0, "Duracell", "19.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Duracell Batterien Optimum AAA/LR03 8 Stück 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905479", "Coop Alkaline Batterie, LR6/AA, 4 Stk.", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/coop-alkaline-batterie-lr6aa-4-stk/p/4905479", "4ST", 7, 4, "Coop", "4.95", "1.24/1ST", "Preis pro 1 Stück", "1.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Coop Alkaline Batterie, LR6/AA, 4 Stk. 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905467", "Alkaline Batterie AAA/LR03 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/alkaline-batterie-aaalr03-4-stueck/p/4905467", "4ST", 4, 4.5, "Coop", "4.95", "1.24/1ST", "Preis pro 1 Stück", "1.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Alkaline Batterie AAA/LR03 4 Stück 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6801782", "Duracell Batterien Plus AAA/LR03 12 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-12-stueck/p/6801782", "12ST", $null, 0, "Duracell", "14.90", "1.24/1ST", "Preis pro 1 Stück", "1.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Duracell Batterien Plus AAA/LR03 12 Stück 50% Aktion 14.90 Schweizer Franken statt 29.85 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905482", "Alkaline Batterie LR14/C 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr14c-2-stueck/p/4905482", "2ST", 3, 4.5, "Coop", "4.95", "2.48/1ST", "Preis pro 1 Stück", "2.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Alkaline Batterie LR14/C 2 Stück 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6153846", "Duracell Knopfzelle CR2032 3V 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2032-3v-2-stueck/p/6153846", "2ST", 1, 5, "Duracell", "9.95", "4.98/1ST", "Preis pro 1 Stück", "4.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Knopfzelle CR2032 3V 2 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905485", "Alkaline Batterie 6LP3146/9V", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-6lp31469v/p/4905485", "1ST", 2, 4.5, "Coop", "4.95", "4.95/1ST", "Preis pro 1 Stück", "4.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Alkaline Batterie 6LP3146/9V 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4949707", "Varta Knopfzellen CR2025 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2025-2-stueck/p/4949707", "2ST", $null, 0, "Varta", "8.95", "4.48/1ST", "Preis pro 1 Stück", "4.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Knopfzellen CR2025 2 Stück 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6761134", "Duracell Batterien PLUS D/LR20 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-dlr20-2-stueck/p/6761134", "2ST", $null, 0, "Duracell", "9.95", "4.98/1ST", "Preis pro 1 Stück", "4.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Batterien PLUS D/LR20 2 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6753556", "Duracell Batterien PLUS AA/LR6 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-plus-aalr6-4-stueck/p/6753556", "4ST", $null, 0, "Duracell", "9.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Duracell Batterien PLUS AA/LR6 4 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6761133", "Duracell Batterien PLUS C/LR14 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-clr14-2-stueck/p/6761133", "2ST", $null, 0, "Duracell", "9.95", "4.98/1ST", "Preis pro 1 Stück", "4.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Batterien PLUS C/LR14 2 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494066", "Varta Longlife Power 9V 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-9v-1er-bli/p/3494066", "1ST", 2, 4.5, "Varta", "8.95", "8.95/1ST", "Preis pro 1 Stück", "8.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Power 9V 1er Bli 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494865", "Varta Electronics CR2430 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2430-1er-bli/p/3494865", "1ST", 1, 5, "Varta", "4.95", "4.95/1ST", "Preis pro 1 Stück", "4.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics CR2430 1er Bli 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6986541", "Duracell Batterie (CR2032, 4 Stück)", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-cr2032-4-stueck/p/6986541", "4ST", $null, 0, "Duracell", "8.95", "2.24/1ST", "Preis pro 1 Stück", "2.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Batterie (CR2032, 4 Stück) 43% Aktion 8.95 Schweizer Franken statt 15.90 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494234", "Varta Electronics V13GA 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13ga-1er-bli/p/3494234", "1ST", $null, 0, "Varta", "2.95", "2.95/1ST", "Preis pro 1 Stück", "2.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics V13GA 1er Bli 2.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6753557", "Duracell Batterien PLUS AAA/LR03 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-4-stueck/p/6753557", "4ST", $null, 0, "Duracell", "9.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Duracell Batterien PLUS AAA/LR03 4 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494233", "Varta Electronics CR2032 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2032-1er-bli/p/3494233", "1ST", 2, 4.5, "Varta", "4.95", "4.95/1ST", "Preis pro 1 Stück", "4.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics CR2032 1er Bli 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6577801", "Duracell Knopfzelle CR2025 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2025-2-stueck/p/6577801", "2ST", $null, 0, "Duracell", "9.95", "4.98/1ST", "Preis pro 1 Stück", "4.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Knopfzelle CR2025 2 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6761135", "Duracell Batterie PLUS 9V/6LR61 1 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-plus-9v6lr61-1-stueck/p/6761135", "1ST", $null, 0, "Duracell", "9.95", "9.95/1ST", "Preis pro 1 Stück", "9.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Duracell Batterie PLUS 9V/6LR61 1 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4014527", "Varta Longlife Batterien AA/LR6 10 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-batterien-aalr6-10-stueck/p/4014527", "10ST", 4, 3, "Varta", "15.95", "1.60/1ST", "Preis pro 1 Stück", "1.60", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife Batterien AA/LR6 10 Stück 15.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6753554", "Duracell Batterien Optimum AA/LR6 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-4-stueck/p/6753554", "4ST", 1, 5, "Duracell", "11.95", "2.99/1ST", "Preis pro 1 Stück", "2.99", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Duracell Batterien Optimum AA/LR6 4 Stück 11.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905484", "Alkaline Batterie LR20/D 2 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr20d-2-stueck/p/4905484", "2ST", 3, 5, "Coop", "5.95", "2.98/1ST", "Preis pro 1 Stück", "2.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Alkaline Batterie LR20/D 2 Stück 5.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494130", "Varta Longlife Power AA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-4er-bli/p/3494130", "4ST", 1, 4, "Varta", "8.95", "2.24/1ST", "Preis pro 1 Stück", "2.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife Power AA 4er Bli 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494133", "Varta Longlife Power AA 8er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-8er-bli/p/3494133", "8ST", 3, 3, "Varta", "14.95", "1.87/1ST", "Preis pro 1 Stück", "1.87", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife Power AA 8er Bli 14.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494909", "Varta Electronics V23GA 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v23ga-1er-bli/p/3494909", "1ST", 1, 5, "Varta", "3.95", "3.95/1ST", "Preis pro 1 Stück", "3.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics V23GA 1er Bli 3.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4145012", "LED Teelichter warmweiss 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/led-teelichter-warmweiss-4-stueck/p/4145012", $null, $null, 0, "Coop", "8.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'beleuchtung', 'taschenlampen-lichter']", "LED Teelichter warmweiss 4 Stück 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4942597", "Skross Adapter Europa-Schweiz", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-europa-schweiz/p/4942597", $null, $null, 0, "Skross", "14.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Skross Adapter Europa-Schweiz 14.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5683906", "Prix Garantie Wasserkocher 1l", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/prix-garantie-wasserkocher-1l/p/5683906", "1ST", 5, 4, "Coop", "19.95", "19.95/1ST", "Preis pro 1 Stück", "19.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "Prix Garantie Wasserkocher 1l 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494065", "Varta Longlife Max Power AA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-max-power-aa-4er-bli/p/3494065", "4ST", 1, 5, "Varta", "9.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife Max Power AA 4er Bli 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6753555", "Duracell Batterien Optimum AAA/LR03 4 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-optimum-aaalr03-4-stueck/p/6753555", "4ST", $null, 0, "Duracell", "11.95", "2.99/1ST", "Preis pro 1 Stück", "2.99", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Duracell Batterien Optimum AAA/LR03 4 Stück 11.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494138", "Varta Longlife Power Batterien AAA/LR03 8 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-batterien-aaalr03-8-stueck/p/3494138", "8ST", 1, 4, "Varta", "14.95", "1.87/1ST", "Preis pro 1 Stück", "1.87", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Varta Longlife Power Batterien AAA/LR03 8 Stück 14.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494131", "Varta Longlife Power AAA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-aaa-4er-bli/p/3494131", "4ST", 3, 3.5, "Varta", "8.95", "2.24/1ST", "Preis pro 1 Stück", "2.24", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Varta Longlife Power AAA 4er Bli 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4358323", "Rayovac Hörgerätebatterien 312 6 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-312-6-stueck/p/4358323", "6ST", 1, 4, "Rayovac", "9.95", "1.66/1ST", "Preis pro 1 Stück", "1.66", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Rayovac Hörgerätebatterien 312 6 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5751576", "satrap Toasty 1 Toaster", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-toasty-1-toaster/p/5751576", $null, 1, 5, "satrap", "29.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "satrap Toasty 1 Toaster 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5763068", "satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-aspira-sine-a700-beutelloser-zyklonstaubsauger/p/5763068", $null, 3, 3.5, "satrap", "119.00", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger 119.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494063", "Varta Longlife Power C 2er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-c-2er-bli/p/3494063", "2ST", $null, 0, "Varta", "7.95", "3.98/1ST", "Preis pro 1 Stück", "3.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Power C 2er Bli 7.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6298048", "Electrolux E201S Classic Staubbeutel", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/electrolux-e201s-classic-staubbeutel/p/6298048", $null, $null, 0, "Electrolux", "9.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Electrolux E201S Classic Staubbeutel 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6973029", "Severin Tischgrill PG 8565", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-tischgrill-pg-8565/p/6973029", $null, 1, 3, "Severin", "49.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "Severin Tischgrill PG 8565 50% Aktion 49.95 Schweizer Franken statt 99.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6266870", "Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück", "/de/haushalt-tier/elektroartikel-batterien/reer-steckdosen-sicherung-kappe-schweiz-8-stueck/p/6266870", $null, $null, 0, "reer", "3.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien']", "Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück 3.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6119284", "satrap Aqua SA10 Wasserkocher", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-aqua-sa10-wasserkocher/p/6119284", $null, 3, 2.5, "satrap", "49.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "satrap Aqua SA10 Wasserkocher 49.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5750424", "satrap Vento HT2 Reisehaartrockner", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/satrap-vento-ht2-reisehaartrockner/p/5750424", $null, 1, 5, "satrap", "29.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'beautygeraete']", "satrap Vento HT2 Reisehaartrockner 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4096751", "Varta Longlife Power Batterien AA/LR6 6 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-batterien-aalr6-6-stueck/p/4096751", "6ST", $null, 0, "Varta", "12.95", "2.16/1ST", "Preis pro 1 Stück", "2.16", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife Power Batterien AA/LR6 6 Stück 12.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494007", "Varta Longlife Power D 2er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-d-2er-bli/p/3494007", "2ST", $null, 0, "Varta", "8.95", "4.48/1ST", "Preis pro 1 Stück", "4.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Power D 2er Bli 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5872158", "satrap Tischventilator Venti 1", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-1/p/5872158", $null, $null, 0, "satrap", "29.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "satrap Tischventilator Venti 1 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4942594", "Skross Adapter World-Schweiz", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-world-schweiz/p/4942594", $null, $null, 0, "Skross", "19.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Skross Adapter World-Schweiz 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4096628", "Varta Longlife Power Batterien AAA/LR03 6 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-batterien-aaalr03-6-stueck/p/4096628", "6ST", $null, 0, "Varta", "12.95", "2.16/1ST", "Preis pro 1 Stück", "2.16", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Power Batterien AAA/LR03 6 Stück 12.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6378875", "Wilkinson Intuition 4in1 Finish Styler", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/wilkinson-intuition-4in1-finish-styler/p/6378875", "1ST", $null, 0, "Wilkinson", "29.95", "29.95/1ST", "Preis pro 1 Stück", "29.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'beautygeraete']", "Wilkinson Intuition 4in1 Finish Styler 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4905486", "Alkaline Batterie 3LR12/4.5V", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-3lr1245v/p/4905486", "1ST", 1, 1, "Coop", "5.95", "5.95/1ST", "Preis pro 1 Stück", "5.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Alkaline Batterie 3LR12/4.5V 5.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5872164", "Satrap Tischventilator Venti 2", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-2/p/5872164", $null, $null, 0, "satrap", "34.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Satrap Tischventilator Venti 2 34.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6579165", "satrap Trimm Set XA 5-in-1 Haar/Trimm-Set", "/de/kosmetik-gesundheit/herrenpflege-rasur/elektrische-herrenrasierer/satrap-trimm-set-xa-5-in-1-haartrimm-set/p/6579165", $null, $null, 0, "satrap", "59.95", $null, $null, $null, $null, "['kosmetik-gesundheit', 'herrenpflege-rasur', 'elektrische-herrenrasierer']", "satrap Trimm Set XA 5-in-1 Haar/Trimm-Set 59.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6007534", "Trend USB-Stick 8 GB", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-8-gb/p/6007534", $null, 1, 5, "Trend", "9.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend USB-Stick 8 GB 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494230", "Varta Electronics V13GS / V357 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13gs--v357-1er-bli/p/3494230", "1ST", $null, 0, "Varta", "4.95", "4.95/1ST", "Preis pro 1 Stück", "4.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics V13GS / V357 1er Bli 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4881282", "Miele FJM HyClean 3D (4 Stück)", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/miele-fjm-hyclean-3d-4-stueck/p/4881282", $null, $null, 0, "Miele", "8.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Miele FJM HyClean 3D (4 Stück) 50% Aktion 8.95 Schweizer Franken statt 17.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494231", "Varta Electronics CR2025 1er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2025-1er-bli/p/3494231", "1ST", $null, 0, "Varta", "4.95", "4.95/1ST", "Preis pro 1 Stück", "4.95", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Electronics CR2025 1er Bli 4.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4589934", "Varta Longlife AA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-aa-4er-bli/p/4589934", "4ST", $null, 0, "Varta", "7.95", "1.99/1ST", "Preis pro 1 Stück", "1.99", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Longlife AA 4er Bli 7.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6831289", "Rowenta Staubsauger RO2933", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/rowenta-staubsauger-ro2933/p/6831289", $null, $null, 0, "Rowenta", "99.90", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Rowenta Staubsauger RO2933 40% Aktion 99.90 Schweizer Franken statt 169.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6753975", "Severin Standgrill mit Grillplatte PG 8563", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-standgrill-mit-grillplatte-pg-8563/p/6753975", $null, $null, 0, "Severin", "74.50", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "Severin Standgrill mit Grillplatte PG 8563 50% Aktion 74.50 Schweizer Franken statt 149.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4119046", "Varta Ultra Lithium AA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-ultra-lithium-aa-4er-bli/p/4119046", "4ST", 1, 5, "Varta", "14.95", "3.74/1ST", "Preis pro 1 Stück", "3.74", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']", "Varta Ultra Lithium AA 4er Bli 14.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4358322", "Rayovac Hörgerätebatterien 13 6 Stück", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-13-6-stueck/p/4358322", "6ST", 1, 5, "Rayovac", "9.95", "1.66/1ST", "Preis pro 1 Stück", "1.66", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Rayovac Hörgerätebatterien 13 6 Stück 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3591269", "Varta Longlife Max Power C 2er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-c-2er-bli/p/3591269", "2ST", 1, 5, "Varta", "8.95", "4.48/1ST", "Preis pro 1 Stück", "4.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Max Power C 2er Bli 8.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6735643", "LED 31V Anschlussset Transf.+Verl.kabel", "/de/haushalt-tier/haushalt-kueche/uebrige-haushaltsartikel/led-31v-anschlussset-transfverlkabel/p/6735643", $null, 1, 5, "Coop", "9.95", $null, $null, $null, $null, "['haushalt-tier', 'haushalt-kueche', 'uebrige-haushaltsartikel']", "LED 31V Anschlussset Transf.+Verl.kabel 50% Aktion 9.95 Schweizer Franken statt 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6700107", "MIELE Compact C2 Classic (mit Beutel, 890 W)", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/miele-compact-c2-classic-mit-beutel-890-w/p/6700107", $null, $null, 0, "Miele", "229.00", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "MIELE Compact C2 Classic (mit Beutel, 890 W) 41% Aktion 229.00 Schweizer Franken statt 389.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5882124", "Philips Avent Audio Monitors DECT-Babyphone", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/philips-avent-audio-monitors-dect-babyphone/p/5882124", $null, $null, 0, "Avent", "99.90", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete']", "Philips Avent Audio Monitors DECT-Babyphone 99.90 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6670192", "Philips Wasserkocher HD9318/01", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/philips-wasserkocher-hd931801/p/6670192", $null, $null, 0, "Philips", "49.90", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "Philips Wasserkocher HD9318/01 49.90 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6288875", "satrap Bügeleisen mit Power-Dampfstation Station XA 6.5", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-buegeleisen-mit-power-dampfstation-station-xa-65/p/6288875", $null, 1, 2, "satrap", "99.90", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "satrap Bügeleisen mit Power-Dampfstation Station XA 6.5 32% Aktion 99.90 Schweizer Franken statt 149.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6425996", "satrap espresso XA Kolbenkaffeemaschine", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-espresso-xa-kolbenkaffeemaschine/p/6425996", $null, 7, 4, "satrap", "249.00", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "satrap espresso XA Kolbenkaffeemaschine 249.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5751209", "satrap Ferru BE1 2200 Dampfbügeleisen", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-ferru-be1-2200-dampfbuegeleisen/p/5751209", $null, 3, 2.5, "satrap", "49.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "satrap Ferru BE1 2200 Dampfbügeleisen 49.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6125818", "satrap Mano XA Handmixer", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mano-xa-handmixer/p/6125818", $null, $null, 0, "satrap", "49.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']", "satrap Mano XA Handmixer 49.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("7016089", "Tefal Bügeleisen FV4961S0", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/tefal-buegeleisen-fv4961s0/p/7016089", $null, $null, 0, "Tefal", "69.90", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Tefal Bügeleisen FV4961S0 30% Aktion 69.90 Schweizer Franken statt 99.90 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6821480", "Tefal Dampfbügeleisen Express Protect SV9202S0 (7.5 Bar, Durilium)", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/tefal-dampfbuegeleisen-express-protect-sv9202s0-75-bar-durilium/p/6821480", $null, $null, 0, "Tefal", "199.00", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Tefal Dampfbügeleisen Express Protect SV9202S0 (7.5 Bar, Durilium) 20% Aktion 199.00 Schweizer Franken statt 249.00 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6822783", "Trend Car Charger QC", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-car-charger-qc/p/6822783", $null, $null, 0, "Trend", "19.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend Car Charger QC 19.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6007537", "Trend USB-Stick 16 GB", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-16-gb/p/6007537", $null, 1, 2, "Trend", "16.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend USB-Stick 16 GB 16.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6459240", "Trend USB-Stick 256 GB", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-256-gb/p/6459240", $null, $null, 0, "Trend", "59.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend USB-Stick 256 GB 59.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("6007538", "Trend USB-Stick 32 GB", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-32-gb/p/6007538", $null, $null, 0, "Trend", "29.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend USB-Stick 32 GB 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5867973", "Trend USB-Stick 64 GB", "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-64-gb/p/5867973", $null, 4, 2, "Trend", "39.95", $null, $null, $null, $null, "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']", "Trend USB-Stick 64 GB 39.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("5894674", "Trisa Beauty Sonic Nail Care System", "/de/kosmetik-gesundheit/make-up/naegel/set-accessoires/trisa-beauty-sonic-nail-care-system/p/5894674", "1ST", $null, 0, "Trisa", "29.95", "29.95/1ST", "Preis pro 1 Stück", "29.95", "1ST", "['kosmetik-gesundheit', 'make-up', 'naegel', 'set-accessoires']", "Trisa Beauty Sonic Nail Care System 29.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4589933", "Varta Longlife AAA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-aaa-4er-bli/p/4589933", "4ST", $null, 0, "Varta", "7.95", "1.99/1ST", "Preis pro 1 Stück", "1.99", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Varta Longlife AAA 4er Bli 7.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("4589935", "Varta Longlife C 2er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-c-2er-bli/p/4589935", "2ST", $null, 0, "Varta", "6.95", "3.48/1ST", "Preis pro 1 Stück", "3.48", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife C 2er Bli 6.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3494067", "Varta Longlife Max Power AAA 4er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-max-power-aaa-4er-bli/p/3494067", "4ST", $null, 0, "Varta", "9.95", "2.49/1ST", "Preis pro 1 Stück", "2.49", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']", "Varta Longlife Max Power AAA 4er Bli 9.95 Schweizer Franken", "2022-08-21 20:58:06"),
  @("3591272", "Varta Longlife Max Power D 2er Bli", "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-d-2er-bli/p/3591272", "2ST", $null, 0, "Varta", "9.95", "4.98/1ST", "Preis pro 1 Stück", "4.98", "1ST", "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']", "Varta Longlife Max Power D 2er Bli 9.95 Schweizer Franken", "2022-08-21 20:58:06")
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

$ws.Rows(89).Delete()
